$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.465.58'
$ws.Range("E2").Value = '  -8.23%  '

$ws.Range("D3").Value = '2.385.49'
$ws.Range("E3").Value = '  -11.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").Value = '459.56'
$ws.Range("E5").Value = '  -7.86%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.80'
$ws.Range("E6").Value = '  -5.33%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '0.481'
$ws.Range("E8").Value = '  -8.38%  '

$ws.Range("D9").Value = '2.400.75'
$ws.Range("E9").Value = '  -10.81%  '

$ws.Range("D10").Value = '0.0944'
$ws.Range("E10").Value = '  -7.20%  '

$ws.Range("D11").Value = '5.26'
$ws.Range("E11").Value = '  -11.65%  '

$ws.Range("D12").Value = '0.313'
$ws.Range("E12").Value = '  -8.25%  '

$ws.Range("E13").Value = '  -4.62%  '

$ws.Range("D14").Value = '2.809.05'
$ws.Range("E14").Value = '  -11.01%  '

$ws.Range("D15").Value = '53.725.13'
$ws.Range("E15").Value = '  -7.75%  '

$ws.Range("D16").Value = '19.58'
$ws.Range("E16").Value = '  -7.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000130'
$ws.Range("E17").Value = '  -2.36%  '

$ws.Range("D18").Value = '2.416.02'
$ws.Range("E18").Value = '  -10.20%  '

$ws.Range("D19").Value = '4.15'
$ws.Range("E19").Value = '  -10.29%  '

$ws.Range("D20").Value = '305.09'
$ws.Range("E20").Value = '  -9.30%  '

$ws.Range("D21").Value = '9.33'
$ws.Range("E21").Value = '  -13.17%  '

$ws.Range("D22").Value = '0.991'
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("D23").Value = '5.68'
$ws.Range("E23").Value = '  +1.36%  '

$ws.Range("D24").Value = '5.31'
$ws.Range("E24").Value = '  -13.02%  '

$ws.Range("D25").Value = '55.92'
$ws.Range("E25").Value = '  -9.49%  '

$ws.Range("E26").Value = '  +1.26%  '

$ws.Range("D27").Value = '2.533.44'
$ws.Range("E27").Value = '  -9.74%  '

$ws.Range("D28").Value = '0.381'
$ws.Range("E28").Value = '  -8.68%  '

$ws.Range("E29").Value = '  -9.85%  '

$ws.Range("D30").Value = '7.09'
$ws.Range("E30").Value = '  -2.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").Value = '0.0₃0709'
$ws.Range("E32").Value = '  -12.46%  '

$ws.Range("D33").Value = '146.06'
$ws.Range("E33").Value = '  -0.70%  '

$ws.Range("D34").Value = '17.57'
$ws.Range("E34").Value = '  -6.42%  '

$ws.Range("D35").Value = '1.42'
$ws.Range("E35").Value = '  -10.46%  '

$ws.Range("D36").Value = '4.96'
$ws.Range("E36").Value = '  -5.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.50'
$ws.Range("E37").Value = '  -14.47%  '

$ws.Range("D38").Value = '1.05'
$ws.Range("E38").Value = '  -4.91%  '

$ws.Range("D39").Value = '0.785'
$ws.Range("E39").Value = '  -13.81%  '

$ws.Range("D40").Value = '0.995'
$ws.Range("E40").Value = '  -0.23%  '

$ws.Range("D41").Value = '32.93'
$ws.Range("E41").Value = '  -8.08%  '

$ws.Range("D42").Value = '0.589'
$ws.Range("E42").Value = '  +0.04%  '

$ws.Range("D43").Value = '3.24'
$ws.Range("E43").Value = '  -6.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0520'
$ws.Range("E44").Value = '  -4.93%  '

$ws.Range("D45").Value = '10.19'
$ws.Range("E45").Value = '  -1.41%  '

$ws.Range("E46").Value = '  -9.91%  '

$ws.Range("D47").Value = '1.938.31'
$ws.Range("E47").Value = '  -9.48%  '

$ws.Range("D48").Value = '0.0216'
$ws.Range("E48").Value = '  -2.87%  '

$ws.Range("D49").Value = '0.0863'
$ws.Range("E49").Value = '  -1.37%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '16.39'
$ws.Range("E50").Value = '  -11.31%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '4.15'
$ws.Range("E51").Value = '  -9.22%  '
